$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.119.64'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.308.57'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'301.29"
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = "'97.37"
$ws.Range("E6").Value = '  -4.34%  '
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("D10").Value = "'33.70"
$ws.Range("E10").Value = '  -4.45%  '
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").Value = "'49.46"
$ws.Range("E12").Value = '  -4.16%  '
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").Value = "'16.95"
$ws.Range("E14").Value = '  +8.42%  '
$ws.Range("D15").Value = "'6.79"
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '2.658.19'
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = '2.284.07'
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").Value = "'0.805"
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").Value = '42.955.93'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = "'11.63"
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = "'6.02"
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").Value = "'67.19"
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").Value = "'236.94"
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = '  -3.12%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'24.90"
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("E29").Value = '  +3.79%  '
$ws.Range("D30").Value = "'166.37"
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = "'33.99"
$ws.Range("D32").Value = "'9.14"
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = "'4.79"
$ws.Range("E34").Value = '  +6.00%  '
$ws.Range("D35").Value = "'4.97"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = "'16.98"
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = "'0.0698"
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").Value = "'2.83"
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("D40").Value = "'0.101"
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = '  -3.89%  '
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").Value = "'2.34"
$ws.Range("E43").Value = '  -2.88%  '
$ws.Range("D44").Value = '1.973.09'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("D46").Value = "'17.63"
$ws.Range("E46").Value = '  -4.90%  '
$ws.Range("D47").Value = "'9.79"
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("D48").Value = "'2.86"
$ws.Range("E48").Value = '  -3.18%  '
$ws.Range("D49").Value = '2.528.72'
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").Value = "'52.88"
$ws.Range("E50").Value = '  -6.93%  '
$ws.Range("D51").Value = "'4.58"
$ws.Range("E51").Value = '  -6.03%  '

Write-Host "Applied 101 cell updates"
